# Update cryptos list prices and volume percentages (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.285.22"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "'2.658.45"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'609.19"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "'149.51"
$ws.Range("E6").Value = "  +3.87%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").Value = "'0.110"
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").Value = "'0.391"
$ws.Range("E10").Value = "  +8.31%  "
$ws.Range("D11").Value = "'5.64"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "'27.85"
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").Value = "'3.126.83"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "'64.011.80"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "'0.0000148"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "'2.640.83"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "'12.01"
$ws.Range("E18").Value = "  +5.15%  "
$ws.Range("D19").Value = "'4.62"
$ws.Range("E19").Value = "  +4.74%  "
$ws.Range("D20").Value = "'347.69"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").Value = "'6.94"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'5.57"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "'66.46"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").Value = "'1.76"
$ws.Range("E25").Value = "  +13.84%  "
$ws.Range("D26").Value = "'1.72"
$ws.Range("E26").Value = "  +4.87%  "
$ws.Range("D27").Value = "'9.40"
$ws.Range("E27").Value = "  +8.43%  "
$ws.Range("D28").Value = "'563.27"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("D29").Value = "'8.28"
$ws.Range("E29").Value = "  +5.67%  "
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'2.08"
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("D33").Value = "'0.0₃0853"
$ws.Range("E33").Value = "  +5.63%  "
$ws.Range("D34").Value = "'1.78"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "'5.35"
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("D36").Value = "'169.00"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").Value = "'0.408"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'1.95"
$ws.Range("E39").Value = "  +5.18%  "
$ws.Range("D40").Value = "'19.36"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'167.53"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").Value = "'40.35"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "'3.86"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("D45").Value = "'0.0574"
$ws.Range("D46").Value = "'22.02"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").Value = "'0.631"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D50").Value = "'0.0965"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "'19.11"
$ws.Range("E51").Value = "  +2.10%  "

# Rows 48 and 49 swap coin data (VeChain and dogwifhat swap ranking positions) with refreshed values
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.00"
$ws.Range("E48").Value = "  +14.18%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0246"
$ws.Range("E49").Value = "  +2.57%  "
